$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Service Index" values (column F) currently sitting in rows 242-292
# actually belong 72 rows further down (rows 314-364). Move them there,
# leaving the original rows blank.
$src = $ws.Range("F242:F292")
$dst = $ws.Range("F314")
$src.Cut($dst)
